$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Match the formatting of the existing header row (e.g. AC1) for the new headers
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# New header cells for the team record columns
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Team record (W/L/T) for every player row
for ($r = 2; $r -le 40; $r++) {
    $ws.Cells.Item($r, 30).Value = 97
    $ws.Cells.Item($r, 31).Value = 65
    $ws.Cells.Item($r, 32).Value = 0
}
